$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(256).Insert()

$ws.Range("A256").Value = 9
$ws.Range("B256").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C256").Value = "Metropolitana"
$ws.Range("D256").Value = 44559
$ws.Range("E256").Value = 13
$ws.Range("F256").Value = 100112031
$ws.Range("G256").Value = "Poroto verde"
$ws.Range("H256").Value = "Magnum"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 34
$ws.Range("K256").Value = 35000
$ws.Range("L256").Value = 37000
$ws.Range("M256").Value = 36000
$ws.Range("N256").Value = "$/malla 25 kilos"
$ws.Range("O256").Value = "Región Metropolitana"
$ws.Range("P256").Value = 1440
$ws.Range("Q256").Value = 25
$ws.Range("R256").Value = "Hortaliza"

Write-Output "done"
